$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New colour codes for the E19:E22 range (panel B options)
$ws.Range("E19").Value = "#DC267F"
$ws.Range("E20").Value = "#648FFF"
$ws.Range("E21").Value = "#FE6100"
$ws.Range("E22").Value = "#5D5C5E"

# Update the selection to reflect the last edited cell
$ws.Range("E22").Select()
